# Generate Report for Archive
#
# The localization status changed from "Ready for handoff" to
# "In Translation" for the tracked document. Update the status cell on
# every sheet that shows it, and shrink the now-narrower "Status"
# columns to match the new (shorter) text, the way Excel's column
# AutoFit would after the text got shorter.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status value everywhere it appears -------------------
# Overview sheet: the per-language status columns (zh-cn / de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-language detail sheets: the "Status" column
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Narrow the Status column(s) to fit the shorter text --------------
# (equivalent to re-running AutoFit on those columns now that the
# status text is shorter)
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
